$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

# Update values on Sheet1
$ws1.Range("A1").Value = 0
$ws1.Range("D1").Value = 3

# Move selection on Sheet3 (no longer the active tab) to A2
$ws3.Range("A2").Select()

# Make Sheet1 the active sheet/tab with selection at H5
$ws1.Activate()
$ws1.Range("H5").Select()
